$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = North Carolina. Fill in the data that was previously an error/blank.
$ws.Range("B4").Value = Get-Date -Year 2020 -Month 7 -Day 4 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"

$ws.Range("C4").Value = 71654
$ws.Range("D4").Value = 1395
$ws.Range("E4").Value = 11390
$ws.Range("F4").Value = 446
$ws.Range("G4").Value = 22.98
$ws.Range("H4").Value = 33.16

$ws.Range("I4").Value = $false
$ws.Range("J4").Value = $true

$ws.Range("K4").Value = 49566
$ws.Range("L4").Value = 1345

$ws.Range("O4").Value = "Success!"
